# Auto-generated edit script applying the Zalera_Profits profit-recalc update
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3699.2568
$ws.Range("I15").Value = 3699.2568
$ws.Range("K15").Value = 11097.7704
$ws.Range("M15").Value = -10928.7704
$ws.Range("H18").Value = 1282.6666
$ws.Range("I18").Value = 1282.6666
$ws.Range("K18").Value = 1282.6666
$ws.Range("M18").Value = -998.6666
$ws.Range("H33").Value = 281.3
$ws.Range("I33").Value = 316.57144
$ws.Range("J33").Value = 199
$ws.Range("K33").Value = 316.57144
$ws.Range("L33").Value = 199
$ws.Range("M33").Value = -87.57144
$ws.Range("N33").Value = -657
$ws.Range("H43").Value = 6390.8184
$ws.Range("I43").Value = 4962.8945
$ws.Range("K43").Value = 4962.8945
$ws.Range("M43").Value = -4893.8945
$ws.Range("H51").Value = 94888.664
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516
$ws.Range("H92").Value = 1093.4375
$ws.Range("I92").Value = 1093.4375
$ws.Range("K92").Value = 1093.4375
$ws.Range("M92").Value = 154.5625
$ws.Range("H103").Value = 665.26666
$ws.Range("I103").Value = 678.5
$ws.Range("J103").Value = 650.1429000000001
$ws.Range("K103").Value = 2035.5
$ws.Range("L103").Value = 1950.4287
$ws.Range("M103").Value = -1449.5
$ws.Range("N103").Value = -3122.4287
$ws.Range("H132").Value = 1345.6129
$ws.Range("I132").Value = 1025.5
$ws.Range("K132").Value = 3076.5
$ws.Range("M132").Value = -546.5
$ws.Range("H138").Value = 3179.8374
$ws.Range("I138").Value = 2664.36
$ws.Range("J138").Value = 3414.1455
$ws.Range("K138").Value = 7993.08
$ws.Range("L138").Value = 10242.4365
$ws.Range("M138").Value = -2853.08
$ws.Range("N138").Value = -20522.4365

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 7857
$ws.Range("H32").Value = 27842.541
$ws.Range("I32").Value = 30623.887
$ws.Range("K32").Value = 30623.887
$ws.Range("M32").Value = -30336.887
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""
$ws.Range("H135").Value = 135690.14
$ws.Range("J135").Value = 135690.14
$ws.Range("L135").Value = 135690.14
$ws.Range("N135").Value = -145830.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 224886
$ws.Range("I86").Value = 1996.2858
$ws.Range("K86").Value = 1996.2858
$ws.Range("M86").Value = -873.2858000000001
$ws.Range("H89").Value = 224886
$ws.Range("I89").Value = 1996.2858
$ws.Range("K89").Value = 9981.429
$ws.Range("M89").Value = -4365.429
$ws.Range("H105").Value = 50014770
$ws.Range("I105").Value = 83355784
$ws.Range("K105").Value = 83355784
$ws.Range("M105").Value = -83354037
$ws.Range("H128").Value = 4000
$ws.Range("I128").Value = 4000
$ws.Range("K128").Value = 12000
$ws.Range("M128").Value = -9510
$ws.Range("H135").Value = 117535.5
$ws.Range("J135").Value = 117535.5
$ws.Range("L135").Value = 117535.5
$ws.Range("N135").Value = -127675.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 156.73914
$ws.Range("I7").Value = 146.9375
$ws.Range("K7").Value = 146.9375
$ws.Range("M7").Value = -33.9375
$ws.Range("H21").Value = 5000
$ws.Range("J21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("N21").Value = -5470
$ws.Range("H22").Value = 485.33334
$ws.Range("I22").Value = 477.14285
$ws.Range("K22").Value = 477.14285
$ws.Range("M22").Value = -127.14285
$ws.Range("H23").Value = 21000
$ws.Range("H27").Value = 21000
$ws.Range("H129").Value = 58593.332
$ws.Range("J129").Value = 58593.332
$ws.Range("L129").Value = 58593.332
$ws.Range("N129").Value = -68593.33199999999
$ws.Range("H134").Value = 11238.083
$ws.Range("J134").Value = 14136.75
$ws.Range("L134").Value = 42410.25
$ws.Range("N134").Value = -47480.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20843284
$ws.Range("I131").Value = 83334380
$ws.Range("J131").Value = 12915.417
$ws.Range("K131").Value = 250003140
$ws.Range("L131").Value = 38746.251
$ws.Range("M131").Value = -249998100
$ws.Range("N131").Value = -48826.251
$ws.Range("H140").Value = 17858408
$ws.Range("I140").Value = 31250740
$ws.Range("K140").Value = 93752220
$ws.Range("M140").Value = -93747040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 50415
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 50415
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 50415
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = -51533
$ws.Range("H111").Value = 30000
$ws.Range("J111").Value = 30000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -36134
$ws.Range("H117").Value = 38950
$ws.Range("J117").Value = 38950
$ws.Range("L117").Value = 38950
$ws.Range("N117").Value = -45834
$ws.Range("H123").Value = 44326.92
$ws.Range("J123").Value = 62500
$ws.Range("L123").Value = 62500
$ws.Range("N123").Value = -67400
$ws.Range("H129").Value = 21932.334
$ws.Range("J129").Value = 21932.334
$ws.Range("L129").Value = 21932.334
$ws.Range("N129").Value = -31932.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = ""
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H133").Value = 104811.664
$ws.Range("J133").Value = 102218
$ws.Range("L133").Value = 102218
$ws.Range("N133").Value = -107278

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 25028
$ws.Range("J40").Value = 25028
$ws.Range("L40").Value = 25028
$ws.Range("N40").Value = -25326
$ws.Range("H74").Value = 14273
$ws.Range("J74").Value = 14273
$ws.Range("L74").Value = 14273
$ws.Range("N74").Value = -16145
$ws.Range("H77").Value = 14273
$ws.Range("J77").Value = 14273
$ws.Range("L77").Value = 42819
$ws.Range("N77").Value = -52179
$ws.Range("H96").Value = 2498.6667
$ws.Range("I96").Value = 2499
$ws.Range("J96").Value = 2498
$ws.Range("K96").Value = 2499
$ws.Range("L96").Value = 2498
$ws.Range("M96").Value = -1126
$ws.Range("N96").Value = -5244
$ws.Range("H107").Value = 1760.2
$ws.Range("I107").Value = 1760.2
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 5280.6
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -3360.6
$ws.Range("N107").Value = ""
$ws.Range("H127").Value = 40971.6
$ws.Range("J127").Value = 44964.5
$ws.Range("L127").Value = 44964.5
$ws.Range("N127").Value = -54884.5
$ws.Range("H132").Value = 6964.5293
$ws.Range("I132").Value = 3121.1428
$ws.Range("J132").Value = 9654.9
$ws.Range("K132").Value = 9363.428400000001
$ws.Range("L132").Value = 28964.7
$ws.Range("M132").Value = -6833.428400000001
$ws.Range("N132").Value = -34024.7
$ws.Range("H136").Value = 4928749.5
$ws.Range("J136").Value = 4415.8887
$ws.Range("L136").Value = 13247.6661
$ws.Range("N136").Value = -18347.6661

